# New weekly price record arrived for "Vega Monumental Concepción -
# Zapallo italiano". It is inserted as row 41 (in chronological order
# with the rest of the sheet), pushing the previous rows 41-111 down to
# 42-112 and growing the used range from A1:R111 to A1:R112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 41; everything below (old 41..111)
# shifts down to 42..112.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with this week's record.
$ws.Range("A41").Value2 = 11
$ws.Range("B41").Value2 = "Vega Monumental Concepción"
$ws.Range("C41").Value2 = "Bíobío"
$ws.Range("D41").Value2 = 44645
$ws.Range("E41").Value2 = 8
$ws.Range("F41").Value2 = 100112032
$ws.Range("G41").Value2 = "Zapallo italiano"
$ws.Range("H41").Value2 = "Sin especificar"
$ws.Range("I41").Value2 = "Primera"
$ws.Range("J41").Value2 = 220
$ws.Range("K41").Value2 = 11000
$ws.Range("L41").Value2 = 12000
$ws.Range("M41").Value2 = 11545
$ws.Range("N41").Value2 = "$/caja 60 unidades"
$ws.Range("O41").Value2 = "Región Metropolitana"
$ws.Range("P41").Value2 = 192
$ws.Range("Q41").Value2 = 60
$ws.Range("R41").Value2 = "Hortaliza"
